$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 800  # was 797
$ws.Range("F3").Value = 2875  # was 2872
$ws.Range("F8").Value = 612  # was 610
$ws.Range("F9").Value = 292  # was 291
$ws.Range("F11").Value = 11917  # was 11899
$ws.Range("F12").Value = 6749  # was 6744
$ws.Range("F21").Value = 290  # was 289
$ws.Range("F23").Value = 3670  # was 3669
$ws.Range("F27").Value = 186  # was 184
$ws.Range("F32").Value = 41  # was 37
$ws.Range("F33").Value = 315  # was 314
$ws.Range("F34").Value = 5058  # was 5057
$ws.Range("F36").Value = 1265  # was 1262
$ws.Range("F38").Value = 698  # was 697
$ws.Range("F39").Value = 1220  # was 1219

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 87  # was 86
$ws.Range("F12").Value = 3701  # was 3700
$ws.Range("F15").Value = 17  # was 16

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9112  # was 9108
$ws.Range("F3").Value = 517  # was 516
$ws.Range("F4").Value = 1861  # was 1859

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9112  # was 9108
$ws.Range("F3").Value = 517  # was 516
$ws.Range("F4").Value = 1861  # was 1859
$ws.Range("F5").Value = 800  # was 797
$ws.Range("F6").Value = 2875  # was 2872
$ws.Range("F12").Value = 292  # was 291
$ws.Range("F14").Value = 11917  # was 11899
$ws.Range("F15").Value = 6749  # was 6744
$ws.Range("F16").Value = 87  # was 86
$ws.Range("F17").Value = 3701  # was 3700
$ws.Range("F25").Value = 290  # was 289
$ws.Range("F27").Value = 3670  # was 3669
$ws.Range("F30").Value = 186  # was 184
$ws.Range("F38").Value = 316  # was 314
$ws.Range("F39").Value = 1265  # was 1262
$ws.Range("F42").Value = 1220  # was 1219
